$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.913.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.668.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.516'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.26'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.904.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.685.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.87%  '
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.929.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.42%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.451.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.69%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.584'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.29%  '
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.813.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.778'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.77%  '
